# Apply the "water monitor commands" sheet update described by the commit
# "update excel file and create pcb".
#
# Net content changes versus the original sheet:
#   1. Column D held a scratch "done" marker column (hecho / x / x x) that is
#      no longer needed -> clear it out.
#   2. The header in F5 changes from "Tipo de consulta" to "Tipo de operacion".
#   3. A new MQTT topic/body pair ("water_tank/full" -> "{capacity: x}") is
#      inserted as a new row right before the old row 39, pushing every row
#      from the old row 39 onward down by one.
#   4. A new example value "{nube: 15.5}" is appended in column G of the new
#      (now empty) row 49.
#   5. The view is left scrolled/selected near the newly added row (E39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the old "hecho / x / x x" helper column entirely.
$ws.Range("D:D").ClearContents()

# 2. Rename the second header column.
$ws.Range("F5").Value = "Tipo de operacion"

# 3. Insert a new blank row at position 39 (shifts rows 39..46 down to 40..47)
#    and fill it with the new "water_tank/full" command.
$ws.Rows("39:39").Insert()
$ws.Range("E39").Value = "water_tank/full"
$ws.Range("G39").Value = "{capacity: x}"

# 4. Add the new cloud-capacity example on the new last row.
$ws.Range("G49").Value = "{nube: 15.5}"

# 5. Update the window view to match (scrolled near row 19, E39 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E39").Select() | Out-Null

Write-Output "applied water monitor commands update"
